# Update the German instruction texts on the "instructions" sheet following
# reviewer feedback:
#  - A3: "Sicherheitsmann" -> "Sicherheitsperson" in the first sentence of the
#        security-guard cover story paragraph.
#  - A9: reworded the sentence about not needing a reason for the preference,
#        replacing "... moegen, sagen Sie uns einfach Ihr Bauchgefuehl." with
#        "... moegen. Teilen Sie uns einfach mit, was Ihr Bauchgefuehl ist."
#
# The worksheet is protected (no password), so it must be unprotected before
# the cell values can be changed, then re-protected afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("instructions")

$ws.Unprotect()

$newA3 = @"
Stellen Sie sich vor, Sie wären ein Sicherheitsperson der auffällige Aktivitäten in einem Unternehmen überwacht. Ihre Aufgabe erfordert ständige Aufmerksamkeit und schnelle Reaktionen, wenn etwas Verdächtiges passiert.   

Wir untersuchen in unserem Labor Aufmerksamkeit und Reaktionsgeschwindigkeit und in diesem Experiment bitten wir Sie, die Rolle des Sicherheitsmannes zu spielen.  

Genauer werden Sie eine Reihe von Dingen auf dem Computerbildschirm beobachten und so schnell wie möglich reagieren, indem Sie die Leertaste drücken, wenn ein Zielgegenstand auftaucht.  


Drücken Sie die Leertaste, um fortzufahren.
"@

$newA9 = @"
Als Nächstes werden Ihnen 30 Paare aus Ziel- und Füll-Wesen aus der Überwachungsaufgabe gezeigt und wir bitten Sie anzugeben, welches Sie lieber mögen.  

Sie brauchen keinen Grund, um eines lieber als das andere zu mögen. Teilen Sie uns einfach mit, was Ihr Bauchgefühl ist.

Uns interessiert, ob die Angenehmheit oder Unangenehmheit der Wesen die Fähigkeit beeinflusst, sie aufmerksam zu beobachten und schnell auf sie zu reagieren. Daher benötigen wir Ihre Angabe, welches Sie lieber mögen.  

Nicht vergessen: Sie brauchen keinen Grund, um eines lieber als das andere zu mögen, also folgen Sie einfach Ihrem Bauchgefühl. Bitte antworten Sie zügig.


Drücken Sie die Leertaste, um fortzufahren. 
"@

$ws.Range("A3").Value = $newA3
$ws.Range("A9").Value = $newA9

$ws.Protect()
